$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) "Codi" row / "Format" cell, paragraph 1:
#    "Numèric de " / "2" / " dígits"  ->  "Alfanumèric de " / "10" / " dígits"
# ---------------------------------------------------------------------------
$paraFormat1 = $d.Paragraphs.Item(6)
$xmlFormat1 = '<w:p ' + $ns + '><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:lang w:eastAsia="ca-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:lang w:eastAsia="ca-ES"/></w:rPr><w:t xml:space="preserve">Alfanumèric de </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:lang w:eastAsia="ca-ES"/></w:rPr><w:t>10</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:lang w:eastAsia="ca-ES"/></w:rPr><w:t xml:space="preserve"> dígits</w:t></w:r></w:p>'
$paraFormat1.Range.InsertXML($xmlFormat1)

# ---------------------------------------------------------------------------
# 2) "Codi" row / "Format" cell, paragraph 2:
#    " enters (99)"  ->  " enters (" / "9999999999" / ")"
# ---------------------------------------------------------------------------
$paraFormat2 = $d.Paragraphs.Item(7)
$xmlFormat2 = '<w:p ' + $ns + '><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:lang w:eastAsia="ca-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:lang w:eastAsia="ca-ES"/></w:rPr><w:t xml:space="preserve"> enters (</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:lang w:eastAsia="ca-ES"/></w:rPr><w:t>9999999999</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:lang w:eastAsia="ca-ES"/></w:rPr><w:t>)</w:t></w:r></w:p>'
$paraFormat2.Range.InsertXML($xmlFormat2)

# ---------------------------------------------------------------------------
# 3) "Codi" row / "Descripció" cell:
#    "Codi únic de dos dígits que identifica el consell comarcal" + ". "
#    -> single long run describing the new 10-digit coding scheme
# ---------------------------------------------------------------------------
$paraDesc = $d.Paragraphs.Item(8)
$descText = "La codificació és a deu dígits. Els dos primers dígits del codi sempre són un 81 (identificador de consell comarcal o Conselh Generau d'Aran), els tres següents són l'identificador de la comarca; per tant, amb un zero al davant (dígit de la tercera posició) tindrem els mateixos codis que les comarques de les quals són corporació local (posicions quarta i cinquena) —això també val per al Conselh Generau d'Aran—; el sisè dígit és un dígit de control, els tres següents són tres zeros i l'últim és un segon dígit de control. Els tres zeros de les posicions setena, vuitena i novena canvien per altres dígits quan es codifiquen els ens de gestió dependents, adscrits o vinculats als consells comarcals i Conselh Generau d'Aran (organismes autònoms, entitats publiques empresarials i societats mercantils participades íntegrament)."
$xmlDesc = '<w:p ' + $ns + '><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:lang w:eastAsia="ca-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:lang w:eastAsia="ca-ES"/></w:rPr><w:t>' + $descText + '</w:t></w:r></w:p>'
$paraDesc.Range.InsertXML($xmlDesc)

# ---------------------------------------------------------------------------
# 4) "Nom" row / "Descripció" cell, last paragraph: drop the _GoBack bookmark
#    that currently sits there (it belongs at the end of the document now).
# ---------------------------------------------------------------------------
$paraNom = $d.Paragraphs.Item(12)
$xmlNom = '<w:p ' + $ns + '><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:lang w:eastAsia="ca-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:lang w:eastAsia="ca-ES"/></w:rPr><w:t xml:space="preserve">Nom que rep </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:lang w:eastAsia="ca-ES"/></w:rPr><w:t>el consell comarcal</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Times New Roman"/><w:color w:val="000000"/><w:sz w:val="20"/><w:lang w:eastAsia="ca-ES"/></w:rPr><w:t>.</w:t></w:r></w:p>'
$paraNom.Range.InsertXML($xmlNom)

# ---------------------------------------------------------------------------
# 5) Final, empty paragraph after the table: this is where the _GoBack
#    bookmark now belongs.
# ---------------------------------------------------------------------------
$paraLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$xmlLast = '<w:p ' + $ns + '><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$paraLast.Range.InsertXML($xmlLast)

# ---------------------------------------------------------------------------
# 6) Drop the now-duplicated "Codi consell comarcal" row entirely (its
#    content was merged into the "Codi" row above). Row deletion must be
#    the last table-structure edit we perform.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.Rows.Item(4).Delete()

Write-Host "done"
